$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7697287797927856
$ws.Range("B1").Value = 2.239279270172119
$ws.Range("C1").Value = 3.345126867294312
$ws.Range("D1").Value = 3.682670116424561
$ws.Range("E1").Value = 0.9134379625320435
